$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")
$ws.Range("A1").Value = "TEST"
